# Fixed update to excel issue
$wb = $excel.ActiveWorkbook

# --- Rename header labels on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet as the 3rd tab ---
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# Move it after "Monthly Trend" so the tab order is Weekly Quantity, Monthly Trend, PO Forecast
$wsForecast.Move($null, $wb.Worksheets.Item("Monthly Trend"))

# NOTE: the sheet handle is positional, so after a Move it may point at
# whatever sheet now occupies the old slot. Re-resolve it by name.
$wsForecast = $wb.Worksheets.Item("PO Forecast")

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(45298.99999999999, 19, -4.164392399777802, 40.32069488573143),
    @(45305.99999999999, 22, 0.4399833164654967, 43.67047654123253),
    @(45319.99999999999, 29, 5.823884891260069, 50.59189044148275),
    @(45403.99999999999, 68, 47.07692214385161, 88.58454371000538),
    @(45410.99999999999, 72, 48.83364112988922, 93.05175808614764),
    @(45417.99999999999, 75, 53.30408655126044, 95.59624138490462),
    @(45424.99999999999, 78, 56.5210685368905, 100.0395858648231),
    @(45431.99999999999, 82, 61.10912256998255, 101.9789224492561),
    @(45438.99999999999, 85, 65.56994553806356, 104.6630497294613),
    @(45445.99999999999, 88, 66.35698798442081, 109.8677678672912),
    @(45452.99999999999, 92, 69.26228840073243, 113.5146341948395),
    @(45459.99999999999, 95, 72.34526352142696, 116.0217487060317)
)

$row = 2
foreach ($entry in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $entry[0]
    $wsForecast.Cells.Item($row, 2).Value = $entry[1]
    $wsForecast.Cells.Item($row, 3).Value = $entry[2]
    $wsForecast.Cells.Item($row, 4).Value = $entry[3]
    $row = $row + 1
}

# Column A holds date-serials; match the same date number format used for
# the "Order Week" / "Order Month" columns on the other two sheets.
$wsForecast.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
